$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2 through 66
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05)
for ($row = 2; $row -le 66; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value = 45174
    }
}
